$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D look like plain numbers (e.g. "243.33",
# "0.0844", "36.424.84" with dots used as thousands separators) so
# Excel would silently reinterpret them as floating point numbers on
# assignment. Forcing a text number format first keeps the exact
# literal string from the source data. Column E (percent strings like
# "  +0.60%  ") is never numeric so it does not need this treatment.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.424.84"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.940.82"
$ws.Range("E3").Value = "  -0.94%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.33"
$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.605"
$ws.Range("E6").Value = "  -2.13%  "

$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("E9").Value = "  -2.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0844"
$ws.Range("E10").Value = "  -0.90%  "

$ws.Range("E11").Value = "  -0.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.226.58"
$ws.Range("E12").Value = "  -0.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.810"
$ws.Range("E13").Value = "  -2.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.22"
$ws.Range("E14").Value = "  -4.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.45"
$ws.Range("E15").Value = "  +0.58%  "

$ws.Range("E16").Value = "  -3.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.941.30"
$ws.Range("E17").Value = "  -0.94%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.374.20"
$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.26"
$ws.Range("E19").Value = "  -2.13%  "

$ws.Range("E20").Value = "  -2.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.53"
$ws.Range("E21").Value = "  -0.85%  "

$ws.Range("E22").Value = "  -2.41%  "

$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("E24").Value = "  -5.55%  "

$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.18"
$ws.Range("E26").Value = "  -3.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.69"
$ws.Range("E27").Value = "  -2.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.133"
$ws.Range("E28").Value = "  +5.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.20"
$ws.Range("E29").Value = "  -3.10%  "

$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("E31").Value = "  -4.48%  "

$ws.Range("E32").Value = "  -3.28%  "

$ws.Range("E33").Value = "  -3.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.24"
$ws.Range("E34").Value = "  +5.36%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.17"
$ws.Range("E35").Value = "  -3.86%  "

$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("E37").Value = "  -0.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.17"
$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.15"
$ws.Range("E39").Value = "  +9.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0987"
$ws.Range("E40").Value = "  +2.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.89"
$ws.Range("E41").Value = "  +0.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0208"
$ws.Range("E42").Value = "  -0.48%  "

$ws.Range("E43").Value = "  -2.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.95"
$ws.Range("E44").Value = "  +1.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.341.24"
$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("E46").Value = "  -2.84%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.19"
$ws.Range("E47").Value = "  -0.93%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.33"
$ws.Range("E48").Value = "  -2.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.81"
$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.118.15"
$ws.Range("E50").Value = "  -0.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.21"
$ws.Range("E51").Value = "  -2.84%  "
